$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename/extend category columns ---
$ws.Range("B1").Value = "Processing"
$ws.Range("C1").Value = "AI-Synonyms"
$ws.Range("D1").Value = "Product-AI"
$ws.Range("E1").Value = "Business-Process-AI"
$ws.Range("F1").Value = "Data"
$ws.Range("G1").Value = "Adjectives"

# G1 is a brand-new header cell - give it the same look (bold/border/center)
# as the rest of the header row by copying the format from the untouched B1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Data rows: year label in column A plus heuristic counts B..G ---
$data = @(
    @{ Year = "2010"; B = 127; C = 5;  D = 10; E = 8; F = 35;  G = 20 },
    @{ Year = "2011"; B = 144; C = 6;  D = 10; E = 6; F = 36;  G = 26 },
    @{ Year = "2012"; B = 142; C = 5;  D = 10; E = 0; F = 45;  G = 35 },
    @{ Year = "2013"; B = 173; C = 10; D = 15; E = 1; F = 32;  G = 80 },
    @{ Year = "2014"; B = 143; C = 12; D = 6;  E = 0; F = 45;  G = 70 },
    @{ Year = "2015"; B = 63;  C = 6;  D = 5;  E = 0; F = 23;  G = 37 },
    @{ Year = "2016"; B = 69;  C = 7;  D = 4;  E = 0; F = 26;  G = 51 },
    @{ Year = "2017"; B = 74;  C = 11; D = 5;  E = 2; F = 33;  G = 70 },
    @{ Year = "2018"; B = 69;  C = 13; D = 8;  E = 2; F = 36;  G = 84 },
    @{ Year = "2019"; B = 180; C = 73; D = 73; E = 1; F = 245; G = 164 },
    @{ Year = "2020"; B = 56;  C = 4;  D = 1;  E = 1; F = 25;  G = 7 }
)

$row = 2
foreach ($entry in $data) {
    # Force the year label to be stored as text (matching the original
    # sharedString-backed labels) rather than auto-detected as a number.
    $ws.Cells.Item($row, 1).NumberFormat = "@"
    $ws.Cells.Item($row, 1).Value = $entry.Year

    $ws.Cells.Item($row, 2).Value = $entry.B
    $ws.Cells.Item($row, 3).Value = $entry.C
    $ws.Cells.Item($row, 4).Value = $entry.D
    $ws.Cells.Item($row, 5).Value = $entry.E
    $ws.Cells.Item($row, 6).Value = $entry.F
    $ws.Cells.Item($row, 7).Value = $entry.G
    $row++
}

# The NumberFormat="@" trick above leaves a stray text-format style on each
# A-column cell; restore the original bold/border/center label look (same
# style already used by the header row) in one batched format-only paste.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A2:A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false
